# Auto-generated edit script: updates profit-calculation columns (H-N)
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per the scheduled price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 55
$ws.Range("H55").Value = 656.0714
$ws.Range("J55").Value = 709.8182
$ws.Range("L55").Value = 709.8182
$ws.Range("N55").Value = -1137.8182
# Row 88
$ws.Range("H88").Value = 4150.8423
$ws.Range("J88").Value = 4474.9414
$ws.Range("L88").Value = 4474.9414
$ws.Range("N88").Value = -5286.9414
# Row 91
$ws.Range("H91").Value = 4150.8423
$ws.Range("J91").Value = 4474.9414
$ws.Range("L91").Value = 4474.9414
$ws.Range("N91").Value = -7282.9414
# Row 105
$ws.Range("H105").Value = 48399.6
$ws.Range("J105").Value = 48399.6
$ws.Range("L105").Value = 48399.6
$ws.Range("N105").Value = -55387.6
# Row 112
$ws.Range("H112").Value = 2860
$ws.Range("I112").Value = 1899.6666
$ws.Range("J112").Value = 3180.111
$ws.Range("K112").Value = 5698.9998
$ws.Range("L112").Value = 9540.332999999999
$ws.Range("M112").Value = -4590.9998
$ws.Range("N112").Value = -11756.333
# Row 137
$ws.Range("H137").Value = 1945.5
$ws.Range("I137").Value = 1530.6
$ws.Range("K137").Value = 4591.799999999999
$ws.Range("M137").Value = -2041.799999999999
# Row 138
$ws.Range("H138").Value = 4377.9165
$ws.Range("I138").Value = 3117.0833
$ws.Range("J138").Value = 4798.1943
$ws.Range("K138").Value = 9351.249899999999
$ws.Range("L138").Value = 14394.5829
$ws.Range("M138").Value = -4211.249899999999
$ws.Range("N138").Value = -24674.5829

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 6792.5835
$ws.Range("I61").Value = 6863.364
$ws.Range("K61").Value = 6863.364
$ws.Range("M61").Value = -6651.364
# Row 63
$ws.Range("H63").Value = 6380
$ws.Range("I63").Value = 950
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 950
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = -264
$ws.Range("N63").Value = -11372
# Row 66
$ws.Range("H66").Value = 6380
$ws.Range("I66").Value = 950
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 4750
$ws.Range("L66").Value = 50000
$ws.Range("M66").Value = -1318
$ws.Range("N66").Value = -56864
# Row 74
$ws.Range("H74").Value = 2351.9
$ws.Range("I74").Value = 2351.9
$ws.Range("K74").Value = 2351.9
$ws.Range("M74").Value = -1477.9
# Row 77
$ws.Range("H77").Value = 2351.9
$ws.Range("I77").Value = 2351.9
$ws.Range("K77").Value = 11759.5
$ws.Range("M77").Value = -7391.5
# Row 106
$ws.Range("H106").Value = 38333.332
$ws.Range("J106").Value = 38333.332
$ws.Range("L106").Value = 38333.332
$ws.Range("N106").Value = -40857.332
# Row 132
$ws.Range("H132").Value = 3917.625
$ws.Range("I132").Value = 4084.5715
$ws.Range("K132").Value = 12253.7145
$ws.Range("M132").Value = -9723.7145
# Row 136
$ws.Range("H136").Value = 6792.5835
$ws.Range("I136").Value = 6863.364
$ws.Range("K136").Value = 20590.092
$ws.Range("M136").Value = -18040.092

$ws = $wb.Worksheets.Item("BSM")
# Row 12
$ws.Range("H12").Value = 101.5
$ws.Range("I12").Value = 52.333332
$ws.Range("J12").Value = 249
$ws.Range("K12").Value = 52.333332
$ws.Range("L12").Value = 249
$ws.Range("M12").Value = 115.666668
$ws.Range("N12").Value = -585
# Row 51
$ws.Range("H51").Value = 133999.5
$ws.Range("J51").Value = 133999.5
$ws.Range("L51").Value = 133999.5
$ws.Range("N51").Value = -134981.5
# Row 94
$ws.Range("H94").Value = 1197.4117
$ws.Range("I94").Value = 1210.4
$ws.Range("K94").Value = 1210.4
$ws.Range("M94").Value = -759.4000000000001
# Row 99
$ws.Range("H99").Value = 5050
# Row 119
$ws.Range("H119").Value = 65217.43
$ws.Range("J119").Value = 65217.43
$ws.Range("L119").Value = 65217.43
$ws.Range("N119").Value = -74893.42999999999
# Row 134
$ws.Range("H134").Value = 3427.375
$ws.Range("I134").Value = 3569.8333
$ws.Range("K134").Value = 10709.4999
$ws.Range("M134").Value = -8174.499899999999

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3607
$ws.Range("I31").Value = 1848.6086
$ws.Range("K31").Value = 1848.6086
$ws.Range("M31").Value = -1553.6086
# Row 34
$ws.Range("H34").Value = 3607
$ws.Range("I34").Value = 1848.6086
$ws.Range("K34").Value = 1848.6086
$ws.Range("M34").Value = -1646.6086
# Row 94
$ws.Range("H94").Value = 2786.36
$ws.Range("I94").Value = 1916.5
$ws.Range("J94").Value = 3589.3076
$ws.Range("K94").Value = 1916.5
$ws.Range("L94").Value = 3589.3076
$ws.Range("M94").Value = -1465.5
$ws.Range("N94").Value = -4491.3076
# Row 132
$ws.Range("H132").Value = 2189.1428
$ws.Range("I132").Value = 2057.5652
$ws.Range("J132").Value = 2794.4
$ws.Range("K132").Value = 6172.6956
$ws.Range("L132").Value = 8383.200000000001
$ws.Range("M132").Value = -3642.6956
$ws.Range("N132").Value = -13443.2

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1210.2858
$ws.Range("I5").Value = 618
$ws.Range("K5").Value = 1854
$ws.Range("M5").Value = -1742
# Row 63
$ws.Range("H63").Value = 2500
$ws.Range("I63").Value = 2500
$ws.Range("K63").Value = 7500
$ws.Range("M63").Value = -6751
# Row 66
$ws.Range("H66").Value = 2500
$ws.Range("I66").Value = 2500
$ws.Range("K66").Value = 22500
$ws.Range("M66").Value = -18756
# Row 87
$ws.Range("H87").Value = 14
$ws.Range("I87").Value = 14
$ws.Range("K87").Value = 42
$ws.Range("M87").Value = 1206
# Row 90
$ws.Range("H90").Value = 14
$ws.Range("I90").Value = 14
$ws.Range("K90").Value = 126
$ws.Range("M90").Value = 6114
# Row 131
$ws.Range("H131").Value = 1841.037
$ws.Range("I131").Value = 846.3333
$ws.Range("J131").Value = 2125.238
$ws.Range("K131").Value = 2538.9999
$ws.Range("L131").Value = 6375.714
$ws.Range("M131").Value = 2501.0001
$ws.Range("N131").Value = -16455.714
# Row 135
$ws.Range("H135").Value = 1210.2858
$ws.Range("I135").Value = 618
$ws.Range("K135").Value = 5562
$ws.Range("M135").Value = -3027
# Row 139
$ws.Range("H139").Value = 5512.8125
$ws.Range("J139").Value = 4636
$ws.Range("L139").Value = 13908
$ws.Range("N139").Value = -24188

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 210
$ws.Range("I5").Value = 210
$ws.Range("K5").Value = 210
$ws.Range("M5").Value = -98
# Row 45
$ws.Range("H45").Value = 37500
$ws.Range("J45").Value = 37500
$ws.Range("L45").Value = 37500
$ws.Range("N45").Value = -38618
# Row 48
$ws.Range("H48").Value = 22357.285
$ws.Range("I48").Value = 13750
$ws.Range("J48").Value = 25800.2
$ws.Range("K48").Value = 13750
$ws.Range("L48").Value = 25800.2
$ws.Range("M48").Value = -13265
$ws.Range("N48").Value = -26770.2
# Row 107
$ws.Range("H107").Value = 1587.0834
$ws.Range("J107").Value = 2073.0715
$ws.Range("L107").Value = 2073.0715
$ws.Range("N107").Value = -5913.0715
# Row 132
$ws.Range("H132").Value = 3262.3333
$ws.Range("I132").Value = 3273.25
$ws.Range("J132").Value = 3175
$ws.Range("K132").Value = 9819.75
$ws.Range("L132").Value = 9525
$ws.Range("M132").Value = -7289.75
$ws.Range("N132").Value = -14585

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 5182.5415
$ws.Range("I40").Value = 3954.8333
$ws.Range("K40").Value = 3954.8333
$ws.Range("M40").Value = -3818.8333
# Row 46
$ws.Range("H46").Value = 2879.5881
$ws.Range("I46").Value = 1099.8
$ws.Range("K46").Value = 1099.8
$ws.Range("M46").Value = -911.8
# Row 122
$ws.Range("H122").Value = 9975.739
$ws.Range("I122").Value = 6716.5625
$ws.Range("J122").Value = 17425.285
$ws.Range("K122").Value = 20149.6875
$ws.Range("L122").Value = 52275.855
$ws.Range("M122").Value = -17699.6875
$ws.Range("N122").Value = -57175.855
# Row 132
$ws.Range("H132").Value = 2263.8823
$ws.Range("I132").Value = 1946
$ws.Range("K132").Value = 5838
$ws.Range("M132").Value = -3308

$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 19429
$ws.Range("J46").Value = 19429
$ws.Range("L46").Value = 19429
$ws.Range("N46").Value = -19891
# Row 113
$ws.Range("H113").Value = 748.4
$ws.Range("I113").Value = 468.8
$ws.Range("K113").Value = 1406.4
$ws.Range("M113").Value = 763.5999999999999
# Row 115
$ws.Range("H115").Value = 95450
$ws.Range("J115").Value = 95450
$ws.Range("L115").Value = 95450
$ws.Range("N115").Value = -98584
# Row 122
$ws.Range("H122").Value = 4604.6665
$ws.Range("I122").Value = 1790.2727
$ws.Range("K122").Value = 5370.8181
$ws.Range("M122").Value = -2920.8181
# Row 134
$ws.Range("H134").Value = 19429
$ws.Range("J134").Value = 19429
$ws.Range("L134").Value = 58287
$ws.Range("N134").Value = -63357
# Row 137
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()
